# Remove 6 closed/duplicate INCO map entries (Caso -162, -195, -197, -200,
# -235, -258) from the "INCO" sheet, matching the data refresh performed by
# the automatic map update (mapa_interactivo_INCO.html).
#
# Row numbers below are the ORIGINAL (pre-delete) 1-based sheet rows that
# hold each Caso to remove. They are deleted highest-row-first so that each
# deletion doesn't shift the row number of one still waiting to be removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToDelete = @(32, 28, 23, 20, 18, 14)

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
